$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update AF5 (related_works) with new list of OpenAlex work URLs
$ws.Range("AF5").Value = 'c("https://openalex.org/W2384274503", "https://openalex.org/W2410163895", "https://openalex.org/W2386504342", "https://openalex.org/W2377136501", "https://openalex.org/W2412005057", "https://openalex.org/W2752085855", "https://openalex.org/W2394784562", "https://openalex.org/W4310603196", "https://openalex.org/W2351826232", "https://openalex.org/W80314699")'

# Update row 11 source/journal metadata fields
$ws.Range("F11").Value = "The American Journal of Gastroenterology"
$ws.Range("G11").Value = "https://openalex.org/S66441642"
$ws.Range("H11").Value = "Lippincott Williams & Wilkins"
$ws.Range("I11").Value = "0002-9270"

# Update any_repository_has_fulltext flag (leading apostrophe forces text so
# Excel doesn't auto-coerce the TRUE/FALSE literal into a native boolean)
$ws.Range("V11").Value = "'FALSE"
